$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Cells.Item(2, 5).Value = 68
$ws.Cells.Item(2, 6).Value = 41
$ws.Cells.Item(2, 8).Value = 41
$ws.Cells.Item(5, 5).Value = 76
$ws.Cells.Item(9, 5).Value = 7
$ws.Cells.Item(10, 5).Value = 228
$ws.Cells.Item(10, 6).Value = 101
$ws.Cells.Item(10, 8).Value = 101
$ws.Cells.Item(11, 5).Value = 175
$ws.Cells.Item(11, 6).Value = 94
$ws.Cells.Item(11, 8).Value = 94
$ws.Cells.Item(12, 5).Value = 252
$ws.Cells.Item(13, 5).Value = 82
$ws.Cells.Item(13, 6).Value = 38
$ws.Cells.Item(13, 8).Value = 38
$ws.Cells.Item(14, 5).Value = 73
$ws.Cells.Item(14, 6).Value = 33
$ws.Cells.Item(14, 8).Value = 33
$ws.Cells.Item(15, 5).Value = 100
$ws.Cells.Item(15, 6).Value = 29
$ws.Cells.Item(15, 8).Value = 29
$ws.Cells.Item(16, 5).Value = 104
$ws.Cells.Item(17, 5).Value = 48
$ws.Cells.Item(17, 6).Value = 23
$ws.Cells.Item(17, 8).Value = 23
$ws.Cells.Item(18, 5).Value = 38
$ws.Cells.Item(20, 5).Value = 60
$ws.Cells.Item(21, 5).Value = 80
$ws.Cells.Item(22, 5).Value = 101
$ws.Cells.Item(22, 6).Value = 47
$ws.Cells.Item(22, 8).Value = 47
$ws.Cells.Item(23, 5).Value = 108
$ws.Cells.Item(24, 5).Value = 118
$ws.Cells.Item(24, 6).Value = 55
$ws.Cells.Item(24, 8).Value = 55
$ws.Cells.Item(25, 5).Value = 116
$ws.Cells.Item(25, 6).Value = 50
$ws.Cells.Item(25, 8).Value = 50
$ws.Cells.Item(26, 5).Value = 70
$ws.Cells.Item(26, 6).Value = 36
$ws.Cells.Item(26, 8).Value = 36
$ws.Cells.Item(27, 5).Value = 166
$ws.Cells.Item(27, 6).Value = 83
$ws.Cells.Item(27, 8).Value = 83
$ws.Cells.Item(28, 5).Value = 102
$ws.Cells.Item(28, 6).Value = 30
$ws.Cells.Item(28, 8).Value = 30
$ws.Cells.Item(29, 5).Value = 108
$ws.Cells.Item(30, 5).Value = 122
$ws.Cells.Item(30, 6).Value = 64
$ws.Cells.Item(30, 8).Value = 64
$ws.Cells.Item(32, 5).Value = 112
$ws.Cells.Item(32, 6).Value = 58
$ws.Cells.Item(32, 8).Value = 58
$ws.Cells.Item(33, 5).Value = 151
$ws.Cells.Item(33, 6).Value = 67
$ws.Cells.Item(33, 8).Value = 67
$ws.Cells.Item(34, 5).Value = 117
$ws.Cells.Item(34, 6).Value = 64
$ws.Cells.Item(34, 8).Value = 64
$ws.Cells.Item(35, 5).Value = 83
$ws.Cells.Item(35, 6).Value = 43
$ws.Cells.Item(35, 8).Value = 43
$ws.Cells.Item(37, 5).Value = 82
$ws.Cells.Item(37, 6).Value = 38
$ws.Cells.Item(37, 8).Value = 38
$ws.Cells.Item(38, 5).Value = 55
$ws.Cells.Item(38, 6).Value = 31
$ws.Cells.Item(38, 8).Value = 31
$ws.Cells.Item(39, 5).Value = 119
$ws.Cells.Item(40, 5).Value = 154
$ws.Cells.Item(40, 6).Value = 64
$ws.Cells.Item(40, 8).Value = 64
$ws.Cells.Item(41, 5).Value = 206
$ws.Cells.Item(41, 6).Value = 79
$ws.Cells.Item(41, 8).Value = 79
$ws.Cells.Item(42, 5).Value = 185
$ws.Cells.Item(42, 6).Value = 91
$ws.Cells.Item(42, 8).Value = 91
$ws.Cells.Item(43, 5).Value = 61
$ws.Cells.Item(43, 6).Value = 26
$ws.Cells.Item(43, 8).Value = 26
$ws.Cells.Item(44, 5).Value = 157
$ws.Cells.Item(44, 6).Value = 72
$ws.Cells.Item(44, 8).Value = 72
$ws.Cells.Item(45, 5).Value = 66
$ws.Cells.Item(46, 5).Value = 143
$ws.Cells.Item(47, 5).Value = 240
$ws.Cells.Item(47, 6).Value = 107
$ws.Cells.Item(47, 8).Value = 107
$ws.Cells.Item(48, 5).Value = 116
$ws.Cells.Item(49, 5).Value = 135
$ws.Cells.Item(49, 6).Value = 57
$ws.Cells.Item(49, 8).Value = 57
$ws.Cells.Item(50, 5).Value = 111
$ws.Cells.Item(50, 6).Value = 45
$ws.Cells.Item(50, 8).Value = 45
$ws.Cells.Item(51, 5).Value = 113
$ws.Cells.Item(51, 6).Value = 46
$ws.Cells.Item(51, 8).Value = 46
$ws.Cells.Item(52, 5).Value = 11
